$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Subtraction RuleDivision Rule" (one paragraph, 3 runs with a spell-check
#    wrapped middle run) -> two separate paragraphs: "Subtraction Rule" and
#    "Division Rule".
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Subtraction RuleDivision Rule", $true, $true, $false, $false, $false, $true, 1, $false, "Subtraction Rule^pDivision Rule", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Remove the "_GoBack" bookmark that currently sits after "Division Rule
#    in Terms of Sets".
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3. Delete the whole "Combinatorial Proofs" paragraph.
# ---------------------------------------------------------------------------
$r = $d.Content
if ($r.Find.Execute("Combinatorial Proofs", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)) {
    $r.Expand(4) | Out-Null   # wdParagraph = 4 -> grab the whole paragraph incl. mark
    $r.Delete()
}

# ---------------------------------------------------------------------------
# 4. Re-create a collapsed "_GoBack" bookmark right after "Binomial
#    coefficient" (inside that paragraph, before its paragraph mark).
#    A collapsed Range placed exactly at "paragraph end - 1" confuses
#    Bookmarks.Add, so a 1-character placeholder is inserted, bookmarked,
#    then emptied back out - Word then keeps the (now collapsed) bookmark
#    anchored at that position.
# ---------------------------------------------------------------------------
$r = $d.Content
if ($r.Find.Execute("Binomial coefficient", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)) {
    $r.Collapse(0) | Out-Null        # wdCollapseEnd
    $r.InsertAfter("Z")
    $rz = $d.Range($r.Start, $r.Start + 1)
    $d.Bookmarks.Add("_GoBack", $rz)
    $rz2 = $d.Range($r.Start, $r.Start + 1)
    $rz2.Text = ""
}

# ---------------------------------------------------------------------------
# 5. "Distinguishable boxes" / "Indistinguishable boxes" -> capitalize the
#    "b" to "B" and split the word into three runs: the leading text (with
#    trailing space), "B" and "oxes" - matching the target markup exactly
#    (no direct formatting on the new runs).
# ---------------------------------------------------------------------------
$boxesRunsXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>B</w:t></w:r><w:r><w:t>oxes</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r = $d.Content
if ($r.Find.Execute("Distinguishable boxes", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)) {
    $boxesStart = $r.Start + 16   # length of "Distinguishable "
    $rBoxes = $d.Range($boxesStart, $r.End)
    $rBoxes.InsertXML($boxesRunsXml)
}

$r = $d.Content
if ($r.Find.Execute("Indistinguishable boxes", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)) {
    $boxesStart = $r.Start + 18   # length of "Indistinguishable "
    $rBoxes = $d.Range($boxesStart, $r.End)
    $rBoxes.InsertXML($boxesRunsXml)
}
